# "Updates to RMI data"
# Swap the battery-cost-decline source from the old BNEF "New Energy Outlook
# 2018" citation to the newer MIT / RSC paper, refresh the learning-rate
# figure on the data sheet to the average of the quoted 20%-27% range, and
# drop the now-stale chart picture + date stamp from the About sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("PDiBCpDoC")

# --- About sheet -----------------------------------------------------

# Drop the reproduced chart picture - it illustrated the retired BNEF source.
for ($i = $ws1.Shapes.Count; $i -ge 1; $i--) {
    $ws1.Shapes.Item($i).Delete()
}

# Drop the "last updated" date stamp that used to sit next to the title.
$ws1.Range("C1").Clear()

# New source citation (was Bloomberg New Energy Finance / New Energy Outlook
# 2018 / https://bnef.turtl.co/story/neo2018 / "Chapter 6, Page 2 ..." ).
$ws1.Range("B3").Value = "Massachusetts Institute of Technology"
$ws1.Range("B4").Value = 2021
$ws1.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$ws1.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$ws1.Range("B7").Value = "Abstract"

# The old "Note that the graph only extends to 2030, not 2050" footnote no
# longer applies (the picture is gone) - blank the cell but keep its style.
$ws1.Range("C8").ClearContents()

# New footnote explaining how the figure below was derived.
$ws1.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# --- PDiBCpDoC sheet ---------------------------------------------------

# Replace the hardcoded 0.18 with the average of the 20%-27% range quoted
# in the new source's abstract.
$ws2.Range("B2").Formula = "=AVERAGE(0.2,0.27)"

# --- Restore cursor positions ------------------------------------------
# Set PDiBCpDoC's selection first, then About's last, so "About" remains
# the active/selected tab (matches the workbook's original state).
$ws2.Range("I4").Select() | Out-Null
$ws1.Range("A10").Select() | Out-Null
